$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to Text format so the new values round-trip
# as literal strings (matching the source data export), not auto-parsed numbers/percentages.
$cells = @("D2","E2","D3","E3","D4","E4","D5","E5","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D28","E28","D40","E40","D41","E41","D42","E42","D43","E43","E44","D45","E45","D47","E47")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "286.37"
$ws.Range("E2").Value = "2.39%"
$ws.Range("D3").Value = "28.57"
$ws.Range("E3").Value = "3.98%"
$ws.Range("D4").Value = "5.070"
$ws.Range("E4").Value = "5.03%"
$ws.Range("D5").Value = "0.06645"
$ws.Range("E5").Value = "3.79%"
$ws.Range("E6").Value = "4.57%"
$ws.Range("D7").Value = "3.408"
$ws.Range("E7").Value = "3.01%"
$ws.Range("D8").Value = "1.371"
$ws.Range("E8").Value = "3.17%"
$ws.Range("D9").Value = "0.9393"
$ws.Range("E9").Value = "4.11%"
$ws.Range("D10").Value = "0.1572"
$ws.Range("E10").Value = "2.25%"
$ws.Range("D11").Value = "0.06592"
$ws.Range("E11").Value = "7.47%"
$ws.Range("D12").Value = "0.07573"
$ws.Range("E12").Value = "1.48%"
$ws.Range("D13").Value = "0.02942"
$ws.Range("E13").Value = "0.65%"
$ws.Range("D14").Value = "0.08974"
$ws.Range("E14").Value = "-0.19%"
$ws.Range("D15").Value = "0.001599"
$ws.Range("E15").Value = "0.12%"
$ws.Range("D16").Value = "0.04495"
$ws.Range("E16").Value = "1.89%"
$ws.Range("D17").Value = "0.0006470"
$ws.Range("E17").Value = "0.64%"
$ws.Range("D18").Value = "0.006274"
$ws.Range("E18").Value = "4.13%"
$ws.Range("D19").Value = "3.443"
$ws.Range("E19").Value = "-1.32%"
$ws.Range("E20").Value = "0.90%"
$ws.Range("D21").Value = "0.3216"
$ws.Range("E21").Value = "2.27%"
$ws.Range("D22").Value = "0.1298"
$ws.Range("E22").Value = "-4.11%"
$ws.Range("D23").Value = "4.058"
$ws.Range("E23").Value = "3.90%"
$ws.Range("D24").Value = "0.1552"
$ws.Range("E24").Value = "3.22%"
$ws.Range("D25").Value = "0.001179"
$ws.Range("E25").Value = "0.43%"
$ws.Range("D26").Value = "0.004146"
$ws.Range("E26").Value = "-3.40%"
$ws.Range("D27").Value = "0.0001249"
$ws.Range("E27").Value = "6.02%"
$ws.Range("D28").Value = "0.0001617"
$ws.Range("E28").Value = "-2.34%"
$ws.Range("D40").Value = "0.04195"
$ws.Range("E40").Value = "2.89%"
$ws.Range("D41").Value = "0.006723"
$ws.Range("E41").Value = "2.03%"
$ws.Range("D42").Value = "0.1250"
$ws.Range("E42").Value = "-10.61%"
$ws.Range("D43").Value = "0.002019"
$ws.Range("E43").Value = "-3.27%"
$ws.Range("E44").Value = "11.95%"
$ws.Range("D45").Value = "0.00005588"
$ws.Range("E45").Value = "0.74%"
$ws.Range("D47").Value = "0.01306"
$ws.Range("E47").Value = "-29.29%"
